$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for columns E (5) through T (20), rows 2-9.
$data = @(
    @(3, 1, 1.182581, 3.547743, 0.2630888587664436, 0.2630888587664436, 3, 1, 0.8180823333333334, 2.454247, 0.5115352725808422, 0.5115352725808422, 0.9674486238356668, 8.707037614521001, 0.1345792310820754, 0.1345792310820754),
    @(3, 1, 1.182581, 3.547743, 0.2630888587664436, 0.2630888587664436, 3, 1, 0.7811863333333333, 2.343559, 0.4884647274191578, 0.4884647274191579, 0.9238161152596668, 8.314345037337, 0.1285096276843682, 0.1285096276843682),
    @(3, 1, 2.308568333333333, 6.925705, 0.513587321461293, 0.5135873214612932, 3, 1, 0.8180823333333334, 2.454247, 0.5115352725808422, 0.5115352725808422, 1.888598968792778, 16.997390719135, 0.2627180304777672, 0.2627180304777672),
    @(3, 1, 2.308568333333333, 6.925705, 0.513587321461293, 0.5135873214612932, 3, 1, 0.7811863333333333, 2.343559, 0.4884647274191578, 0.4884647274191579, 1.803422031566111, 16.230798284095, 0.2508692909835259, 0.2508692909835259),
    @(3, 1, 0.7305033333333334, 2.19151, 0.1625151159132014, 0.1625151159132014, 3, 1, 0.8180823333333334, 2.454247, 0.5115352725808422, 0.5115352725808422, 0.5976118714411112, 5.37850684297, 0.08313221411716663, 0.08313221411716663),
    @(3, 1, 0.7305033333333334, 2.19151, 0.1625151159132014, 0.1625151159132014, 3, 1, 0.7811863333333333, 2.343559, 0.4884647274191578, 0.4884647274191579, 0.5706592204544445, 5.13593298409, 0.07938290179603476, 0.07938290179603476),
    @(2, 0.6666666666666666, 0.2733343333333333, 0.820003, 0.06080870385906196, 0.06080870385906197, 3, 1, 0.8180823333333334, 2.454247, 0.5115352725808422, 0.5115352725808422, 0.2236099891934445, 2.012489902741, 0.03110579690383297, 0.03110579690383297),
    @(2, 0.6666666666666666, 0.2733343333333333, 0.820003, 0.06080870385906196, 0.06080870385906197, 3, 1, 0.7811863333333333, 2.343559, 0.4884647274191578, 0.4884647274191579, 0.2135250456307778, 1.921725410677, 0.02970290695522899, 0.029702906955229)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 5
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
